$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "procura pela polícia" table had a pandas multi-index header exported
# with a stray "unnamed: 1_level_1" placeholder in B2, plus two label-only
# sub-header rows ("situação do domicílio" and "grandes regiões e unidades
# da federação") that carried no data. Correct the header text and remove
# the two empty sub-header rows so the data rows below shift up and close
# the gaps (dados corrigidos / inicio da analise PNAD 2009).

# Fix the stray pandas multi-index placeholder in the header row.
$ws.Range("B2").Value = "total"

# Remove the two label-only rows (delete the lower one first so the row
# index of the upper one doesn't shift before it's deleted).
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(5).Delete()
